$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 updates
$ws.Range("H4").Value = 5.5
$ws.Range("I4").Value = 8.75
$ws.Range("J4").Value = 1.57
$ws.Range("K4").Value = 2.8
$ws.Range("L4").Value = 7.1
$ws.Range("N4").Value = 17.5
$ws.Range("O4").Value = 1.08
$ws.Range("P4").Value = 6.5
$ws.Range("S4").Value = 1.18
$ws.Range("T4").Value = 4.3
$ws.Range("U4").Value = 1.71
$ws.Range("V4").Value = 2.08
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 6.8
$ws.Range("Y4").Value = 8
$ws.Range("Z4").Value = 7.2
$ws.Range("AC4").Value = 21
$ws.Range("AD4").Value = 10.75
$ws.Range("AH4").Value = 27
$ws.Range("AJ4").Value = 23
$ws.Range("AK4").Value = 175
$ws.Range("AL4").Value = 70
$ws.Range("AM4").Value = 55
$ws.Range("AN4").Value = 3.4
$ws.Range("AP4").Value = 12.5
$ws.Range("AQ4").Value = 11.25
$ws.Range("AR4").Value = 28
$ws.Range("AS4").Value = 120
$ws.Range("AT4").Value = 4.1
$ws.Range("AV4").Value = 55
$ws.Range("AW4").Value = 10.25
$ws.Range("AX4").Value = 50
$ws.Range("AY4").Value = 37
$ws.Range("AZ4").Value = 300
$ws.Range("BA4").Value = 250
$ws.Range("BB4").Value = 400

# Row 5 updates
$ws.Range("G5").Value = 1.39
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 7.2
$ws.Range("J5").Value = 1.87
$ws.Range("L5").Value = 6.4
$ws.Range("R5").Value = 2.05
$ws.Range("X5").Value = 6
$ws.Range("Y5").Value = 6.8
$ws.Range("AB5").Value = 18.5
$ws.Range("AE5").Value = 13.5
$ws.Range("AF5").Value = 50
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 17.5
$ws.Range("AI5").Value = 40
$ws.Range("AJ5").Value = 17.5
$ws.Range("AL5").Value = 60
$ws.Range("AO5").Value = 6.4
$ws.Range("AP5").Value = 15
$ws.Range("AQ5").Value = 18
$ws.Range("AS5").Value = 175
$ws.Range("AU5").Value = 7.5
$ws.Range("AW5").Value = 8.5
$ws.Range("AX5").Value = 40
$ws.Range("AY5").Value = 37
$ws.Range("AZ5").Value = 300

$wb.Save()
